# 自动更新Excel文件 - 将"剩余"天数按时间推进一天重新计算
# For every data row: E (剩余/remaining) = D (总天/total days) - (today - F (开始时间/start date)) in days.
# When that countdown reaches zero (or below), a new cycle starts: F resets to "today" and E resets to D.
#
# The workbook snapshot implies "today" moved from 2025-12-16 to 2025-12-17 (one day later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = [datetime]::ParseExact("20251217", "yyyyMMdd", $null)
$todayStr = $today.ToString("yyyyMMdd")

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {

    $dVal = $ws.Cells.Item($row, 4).Value2
    $fVal = $ws.Cells.Item($row, 6).Value2

    if ($dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]([long]$fVal)

    # Skip malformed / non 8-digit date values (e.g. data-entry typos) - leave untouched.
    if ($fStr.Length -ne 8) {
        continue
    }

    $startDate = $null
    try {
        $startDate = [datetime]::ParseExact($fStr, "yyyyMMdd", $null)
    } catch {
        continue
    }

    $totalDays = [int]$dVal
    $elapsed = [Math]::Round($today.ToOADate() - $startDate.ToOADate())
    $remaining = $totalDays - [int]$elapsed

    if ($remaining -le 0) {
        # Cycle finished - start a new one today.
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = [int]$todayStr
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining
    }
}
